$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - LED1: Mid/Ref/Pad X & Y moved, rotated 90
$ws.Range("D13").Value = "41.23mm"
$ws.Range("E13").Value = "14.372mm"
$ws.Range("F13").Value = "41.23mm"
$ws.Range("G13").Value = "14.372mm"
$ws.Range("H13").Value = "41.23mm"
$ws.Range("I13").Value = "15.422mm"
$ws.Range("L13").Value = 90

# Row 14 - LED2: Mid/Ref/Pad X & Y moved, rotated 90
$ws.Range("D14").Value = "37.518mm"
$ws.Range("E14").Value = "14.277mm"
$ws.Range("F14").Value = "37.518mm"
$ws.Range("G14").Value = "14.277mm"
$ws.Range("H14").Value = "37.518mm"
$ws.Range("I14").Value = "15.327mm"
$ws.Range("L14").Value = 90

# Row 21 - R5: Mid/Ref/Pad X & Y moved, rotated 90
$ws.Range("D21").Value = "36.924mm"
$ws.Range("E21").Value = "18.69mm"
$ws.Range("F21").Value = "36.924mm"
$ws.Range("G21").Value = "18.69mm"
$ws.Range("H21").Value = "36.924mm"
$ws.Range("I21").Value = "17.69mm"
$ws.Range("L21").Value = 90

# Row 22 - R6: Mid/Ref/Pad X & Y moved, rotation back to 0
$ws.Range("D22").Value = "40.693mm"
$ws.Range("E22").Value = "17.96mm"
$ws.Range("F22").Value = "40.693mm"
$ws.Range("G22").Value = "17.96mm"
$ws.Range("H22").Value = "39.693mm"
$ws.Range("I22").Value = "17.96mm"
$ws.Range("L22").Value = 0

# Row 27 - RGB: Mid/Ref/Pad X moved (Y unchanged)
$ws.Range("D27").Value = "22.987mm"
$ws.Range("F27").Value = "22.987mm"
$ws.Range("H27").Value = "22.987mm"

# Row 30 - U2: Mid/Ref X & Y and Pad X & Y moved
$ws.Range("D30").Value = "22.606mm"
$ws.Range("E30").Value = "24.892mm"
$ws.Range("F30").Value = "22.606mm"
$ws.Range("G30").Value = "24.892mm"
$ws.Range("H30").Value = "25.356mm"
$ws.Range("I30").Value = "20.642mm"
